# report 年度指标增加 alpha 和 beta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("年度指标")

# remember which sheet was active so we can restore it at the end
$originalActive = $excel.ActiveSheet

# --- new column headers (row 1, human readable labels) ---
$ws.Range("Q1").Value = "阿尔法"
$ws.Range("R1").Value = "贝塔"

# --- new column template placeholders (row 2) ---
$ws.Range("Q2").Value = "#alpha#"
$ws.Range("R2").Value = "#beta#"

# carry over the same formatting used by the preceding column (P)
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null

$ws.Range("P2").Copy() | Out-Null
$ws.Range("Q2:R2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# match the column widths of the other data columns
$ws.Columns.Item(17).ColumnWidth = 23.86
$ws.Columns.Item(18).ColumnWidth = 23.86

# update the sheet selection to the new bottom-right area
$ws.Activate()
$ws.Range("R7").Select() | Out-Null

# restore the originally active sheet/tab
$originalActive.Activate()
